$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Text = $old
    $rng.Find.MatchCase = $true
    $rng.Find.MatchWholeWord = $false
    $rng.Find.MatchWildcards = $false
    $rng.Find.Execute() | Out-Null
    if (-not $rng.Find.Found) {
        throw "Text not found: $old"
    }
    $rng.Text = $new
}

Replace-Text "Who do I contact if I have questions or concerns?" "Wie kan ek kontak as ek vrae of bekommernisse het?"
Replace-Text "If you have any questions or concerns about your rights as a study participant, you can contact the study team at " "As jy enige vrae of bekommernisse het oor jou regte as 'n studie-deelnemer, kan jy die studiespan kontak by "
Replace-Text " or on WhatsApp at +27 XX XXX XXXX (messages only)." " of via WhatsApp by +27 XX XXX XXXX (net boodskappe)."
Replace-Text "If you have more questions or concerns about your rights, you can contact one of the ethics committees listed: " "As jy meer vrae of bekommernisse het oor jou regte, kan jy een van die etiekkomitees hieronder kontak: "
Replace-Text "Name" "Naam"
Replace-Text "Telephone" "Telefoon"
Replace-Text "Email" "E-pos"
Replace-Text "University of Cape Town Centre for Social Science Research" "Universiteit van Kaapstad Sentrum vir Sosiale Wetenskap Navorsing"
Replace-Text "Human Research Ethics Committee" "Etiekkomitee vir Menslike Navorsing"
Replace-Text "Informed Consent to Take Part in the Study" "Ingeligte Toestemming om aan die Studie Deel te Neem"
Replace-Text "Please read these statements carefully:" "Lees asseblief hierdie stellings sorgvuldig:"
Replace-Text "I have read the information above and know what is expected of my child." "Ek het die inligting hierbo gelees en weet wat van my kind verwag word."
Replace-Text "I understand as my child’s guardian that I am giving consent for them to participate. " "Ek verstaan as my kind se voog dat ek toestemming gee vir hulle om deel te neem. "
Replace-Text "I understand that even though I have given consent that my child will still be able to choose freely if they want to be interviewed. " "Ek verstaan dat, alhoewel ek toestemming gegee het, my kind steeds vry sal wees om te kies of hulle ondervra wil word. "
Replace-Text "I understand that they can say no to being interviewed without any consequence. " "Ek verstaan dat hulle nee kan sê om ondervra te word sonder enige gevolge. "
Replace-Text "I know who can see my child’s information, how it will be kept safe, and what happens to it after the study." "Ek weet wie my kind se inligting kan sien, hoe dit veilig gehou sal word, en wat daarmee sal gebeur na die studie."
Replace-Text "I understand that I will not be notified of my child’s answers." "Ek verstaan dat ek nie van my kind se antwoorde in kennis gestel sal word nie."
Replace-Text "I know I and my child won’t be named in any papers or reports from this study." "Ek weet dat ek en my kind nie in enige artikels of verslae van hierdie studie sal genoem word nie."
Replace-Text "I know who to contact if I have a problem with the study." "Ek weet wie ek kan kontak as ek 'n probleem met die studie het."
Replace-Text "You can contact me again if more information is needed from me." "Jy kan my weer kontak as daar meer inligting van my benodig word."
Replace-Text "You can keep my contact information safe so you can tell me about the results of the study." "Jy kan my kontakbesonderhede veilig hou sodat jy my oor die resultate van die studie kan inlig."

Write-Output "done"
